$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "578.12", "63.347.00").
# Force text format before assignment so Excel does not auto-convert them to numbers
# and so values like "1.00" keep their trailing zeros, matching the source export format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.347.00'
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.483.31'
$ws.Range("E3").Value = '  +2.91%  '

$ws.Range("E4").Value = '  -0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.12'
$ws.Range("E5").Value = '  +1.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.79'
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.480.54'
$ws.Range("E9").Value = '  +1.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.112'
$ws.Range("E10").Value = '  +1.68%  '

$ws.Range("E11").Value = '  +1.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.31'
$ws.Range("E12").Value = '  +0.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +1.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.24'
$ws.Range("E14").Value = '  +9.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.930.36'
$ws.Range("E16").Value = '  +2.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.304.41'
$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.483.26'
$ws.Range("E18").Value = '  +1.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.02'
$ws.Range("E19").Value = '  +0.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.13'
$ws.Range("E20").Value = '  +2.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.54'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("E22").Value = '  +9.83%  '

$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.14'
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  +0.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.50'
$ws.Range("E25").Value = '  +1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '667.02'
$ws.Range("E26").Value = '  +7.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("E27").Value = '  +12.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000100'
$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.48'
$ws.Range("E31").Value = '  +4.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  +0.79%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("E33").Value = '  +2.67%  '

$ws.Range("E34").Value = '  +0.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").Value = '  +4.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.82'
$ws.Range("E37").Value = '  +1.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.55'
$ws.Range("E38").Value = '  +2.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.375'
$ws.Range("E39").Value = '  +0.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.09'
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.85'
$ws.Range("E41").Value = '  +1.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.76'
$ws.Range("E42").Value = '  +3.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.77'
$ws.Range("E43").Value = '  +1.29%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0302'
$ws.Range("E45").Value = '  +7.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.18'
$ws.Range("E46").Value = '  +27.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.21'
$ws.Range("E47").Value = '  +3.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.65'
$ws.Range("E48").Value = '  +1.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.99'
$ws.Range("E49").Value = '  +3.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.610'
$ws.Range("E50").Value = '  +1.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0517'
$ws.Range("E51").Value = '  +0.74%  '
